$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D5","D6","D7","D8","D13","D14","D15","D16","D20","D22","D26","D27","D28","D30","D32","D34","D36","D38","D39","D40","D42","D43","D45","D47","D48","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.231.48"
$ws.Range("E2").Value = "  +3.02%  "

$ws.Range("D3").Value = "2.121.54"
$ws.Range("E3").Value = "  +3.43%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "235.08"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("D7").Value = "58.21"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +3.14%  "

$ws.Range("E10").Value = "  +3.87%  "

$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").Value = "2.433.85"
$ws.Range("E12").Value = "  +3.53%  "

$ws.Range("D13").Value = "14.55"
$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("D14").Value = "21.49"
$ws.Range("E14").Value = "  +3.56%  "

$ws.Range("D15").Value = "0.787"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").Value = "5.25"
$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("D17").Value = "2.114.40"
$ws.Range("E17").Value = "  +3.33%  "

$ws.Range("D18").Value = "38.092.79"
$ws.Range("E18").Value = "  +2.87%  "

$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").Value = "70.58"
$ws.Range("E20").Value = "  +2.60%  "

$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").Value = "228.44"
$ws.Range("E22").Value = "  +1.96%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("D26").Value = "169.00"
$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("D27").Value = "0.140"
$ws.Range("E27").Value = "  +11.96%  "

$ws.Range("D28").Value = "9.00"
$ws.Range("E28").Value = "  +2.98%  "

$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").Value = "19.59"
$ws.Range("E30").Value = "  +3.21%  "

$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("D32").Value = "4.65"
$ws.Range("E32").Value = "  +4.97%  "

$ws.Range("E33").Value = "  +3.78%  "

$ws.Range("D34").Value = "0.0626"
$ws.Range("E34").Value = "  +2.60%  "

$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("D36").Value = "3.48"
$ws.Range("E36").Value = "  +6.87%  "

$ws.Range("E37").Value = "  +5.07%  "

$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  -4.46%  "

$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  +8.24%  "

$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").Value = "97.41"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.465.15"
$ws.Range("E44").Value = "  -1.09%  "

$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("E46").Value = "  +4.91%  "

$ws.Range("D47").Value = "4.12"
$ws.Range("E47").Value = "  -8.47%  "

$ws.Range("D48").Value = "15.75"
$ws.Range("E48").Value = "  +3.59%  "

$ws.Range("D49").Value = "3.05"
$ws.Range("E49").Value = "  +3.99%  "

$ws.Range("D50").Value = "7.30"
$ws.Range("E50").Value = "  +2.82%  "

$ws.Range("D51").Value = "2.318.27"
$ws.Range("E51").Value = "  +3.52%  "
